$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain plain text even though many values look numeric
# (e.g. "327.00", "1.000") - set the cell format to Text before assigning the string
# so Excel does not silently coerce them into numbers and strip formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.108.63"
$ws.Range("E2").Value = "  +5.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.920.19"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.00"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5158"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4011"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08450"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.79"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.122"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.71"
$ws.Range("E12").Value = "  +5.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.348"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.921.15"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.352"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.22"
$ws.Range("E17").Value = "  +5.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001116"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06741"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.09"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.060"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.121.24"
$ws.Range("E23").Value = "  +5.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.141.79"
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.70"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.07"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.457"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.30"
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.075"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1060"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.076"
$ws.Range("E33").Value = "  +4.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.664"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02513"
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06600"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2220"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.236"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.019"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.204"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6548"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.241"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.41"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6136"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.765"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.054"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.64"
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.243"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.157"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.30"
$ws.Range("E51").Value = "  +3.18%  "
